# The edit re-shuffles the 14 data rows (rows 2-15) of the "Artfynd" sheet:
# each row's entire record (columns A-AY) ends up holding the data that
# originally lived in a *different* row of the same block - i.e. a pure
# row permutation, presumably caused by re-joining the export against a
# slightly different row order upstream. The header row (row 1) and the
# set of columns used are unchanged.
#
# Strategy: snapshot the whole A2:AY15 block into a 2-D array, build a
# re-ordered array per the explicit before-row -> after-row mapping derived
# from the diff, then write it back in one shot. A couple of the text
# columns (Y/Z/AA/AB) hold values that look like dates/times (e.g.
# "2023-08-15", "00:00"); Excel's Range.Value setter auto-coerces such
# look-alike strings to real date/time serials unless the target cells are
# already formatted as Text, so those columns are pre-formatted before the
# write and restored to the default style afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 15
$firstCol = 1   # A
$lastCol = 51   # AY

$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$snapshot = $srcRange.Value()

# Columns that contain date-/time-looking text ("2023-08-15", "00:00") and
# therefore need to be protected from Excel's automatic type coercion when
# the values are written back via Range.Value.
$dateLikeRange = $ws.Range($ws.Cells.Item($firstRow, 25), $ws.Cells.Item($lastRow, 28))
$dateLikeRange.NumberFormat = "@"

$rowCount = $lastRow - $firstRow + 1
$colCount = $lastCol - $firstCol + 1

# Destination sheet-row -> source sheet-row (which row's original content
# should end up at the destination), taken from the diff.
$rowMap = @{
    2  = 15
    3  = 9
    4  = 14
    5  = 12
    6  = 8
    7  = 11
    8  = 6
    9  = 13
    10 = 4
    11 = 5
    12 = 3
    13 = 2
    14 = 7
    15 = 10
}

$result = New-Object 'object[,]' $rowCount, $colCount

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $destIdx = $destRow - $firstRow
    $srcIdx = $srcRow - $firstRow
    for ($c = 0; $c -lt $colCount; $c++) {
        $result[$destIdx, $c] = $snapshot[$srcIdx + 1, $c + 1]
    }
}

$srcRange.Value = $result

# Restore the default (General) style on the date/time-like columns now
# that the text values are safely committed, so formatting matches the
# original workbook.
$dateLikeRange.Style = "Normal"
